$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 251, shifting the existing rows (251-267) down to (252-268).
$ws.Rows(251).Insert()

# Match the date cell format (column D) used by the surrounding rows.
$ws.Cells.Item(251, 4).NumberFormat = $ws.Cells.Item(252, 4).NumberFormat

# Fill the new row 251 with the inserted record's values.
$ws.Cells.Item(251, 1).Value = 5
$ws.Cells.Item(251, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(251, 3).Value = "Maule"
$ws.Cells.Item(251, 4).Value = 44585
$ws.Cells.Item(251, 5).Value = 7
$ws.Cells.Item(251, 6).Value = 100114013
$ws.Cells.Item(251, 7).Value = "Zanahoria"
$ws.Cells.Item(251, 8).Value = "Sin especificar"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 300
$ws.Cells.Item(251, 11).Value = 8000
$ws.Cells.Item(251, 12).Value = 8000
$ws.Cells.Item(251, 13).Value = 8000
$ws.Cells.Item(251, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(251, 15).Value = "Región de Ñuble"
$ws.Cells.Item(251, 16).Value = 400
$ws.Cells.Item(251, 17).Value = 20
$ws.Cells.Item(251, 18).Value = "Hortaliza"
